$d = $word.ActiveDocument

# --- Fill the previously empty paragraph (Architecting the Dessau building / Wassily Chair) ---
$p11 = $d.Paragraphs.Item(11)
$r11 = $p11.Range
$r11.InsertAfter("Architecting the Dessau building was a reflection of the Bauhaus’ schools new direction.  Under the champion of teacher Maholy-Nagy, the curriculum turned toward functional art that could serve a modern industrial society.  This meant that in addition to the study of architecture itself, it was also the objects within buildings that the artists focused on.  Fixtures were prime targets to transform using Bauhaus principles into utilitarian but also beautiful works.  The Model B3 chair, also known as the Wassily Chair, was designed by Marcel Breuer and named for his colleague Wassily Kandinsky.  This chair has a tubular metal frame, adorned with strips of fabric or leather to construct a minimalist, exposed chair.   As Eskilson describes it, “Its spare steel frame forms cubic shaped that seem to pass through each other, its beauty resting in proportion and the balance of simple forms”  (217).  Another iconic legacy of the Bauhaus school, “the chair…has been mass-produced since the 1950s” (Hartov).  Even today, a quick search will easily yield millions of results, from YouTube videos touting “How to Identify…Wassily Chairs Authenticity” to the many sellers who offer replicas of the chair.  The impact of this chair is evident in its legacy; a beautiful form that people cherish and sit in today.  ")
$r11.Font.Name = "Times New Roman"
$r11.Font.NameBi = "Times New Roman"

# --- Replace "The emphasis..." paragraph text with the first new paragraph (Graphics and typography) ---
$p12 = $d.Paragraphs.Item(12)
$r12 = $p12.Range
$r12.Text = "Graphics and typography also flourished in the new location of the school, turning from the traditional fine arts instruction to instead focus on commercial application.  Despite blocky forms and solid blocks, the arrangement of typography allowed for kinetic infusion into posters and balance.  Herbert Bayer exemplifies the style of the time with posters that show blocked colors and text, but turned and made perpendicular to each other to draw the eye across the composition.  Sans serif type was highly favored and thought to represent the spirit of the machine age and also served simplicity to complement photography.  Typography, coupled with the increased emphasis on technical execution of photography led to the birth of the new typophoto visualization style.  Working together, typography and photography could create a cohesive, unified message across planes and mediums.  "
$r12.Font.Name = "Times New Roman"
$r12.Font.NameBi = "Times New Roman"

# --- Append the remaining new paragraphs after paragraph 12 ---
$insertRange = $d.Paragraphs.Item(12).Range
# paragraph 1
$insertRange.InsertParagraphAfter()
$insertRange = $d.Paragraphs.Item(13).Range
$insertRange.Text = "One of the most famous typefaces created by the Bauhaus was Herbert Bayer’s Universal.  Universal was characterized by even weight thickness, symmetry, and perfect clarity.  The letter forms themselves were designed with great care, from the reflective “n” and “u” letterforms to standardized angles that strokes adhere to while composing the letterforms. Every aspect was intentional.  Universal was designed to be used with a single case letterset (lowercase) that would save printers money without sacrificing readability.  Universal inspired one of the most influential sans serif fonts of all, Futura.  While Futura was not made in Bauhaus, the original intention was to take pure geometric forms and allow for readability.  However, multiple iterations left Futura deviating from the pure geometry because “the purest geometric forms neither appeared beautiful as individual shapes nor connected fluidly with one another” (Eskilson, 224).  In the end Futura deviates from those pure geometric intentions with some subtle strokes to aid readability, but is still heavily inspired by the same principles espoused by Bauhaus.  Futura is still widely used today, as evidenced by its use by the popular food magazine and empire, bon appetit.  Even today, we as consumers are exposed to Bauhaus-inspired typography commonly, further solidifying its place in our everyday.  "
$insertRange.Font.Name = "Times New Roman"
$insertRange.Font.NameBi = "Times New Roman"
# paragraph 2
$insertRange.InsertParagraphAfter()
$insertRange = $d.Paragraphs.Item(14).Range
$insertRange.Text = "The prominence and success of the Bauhaus school meant that it also drew the eye of a darker movement rising in Germany.  The National Socialist German Workers’ Party, colloquially known as the Nazi Party did not support the design and free thinking style of the Bauhaus.  “The Nazis saw the Bauhaus as representing ‘foreignness’ and viewed their designs as distinctly un-German and criticised their modernist style, so when the party gain control of Dessau city council in 1931, they moved to close the school” (“100 Years of Bauhaus”).  The director at the time, Ludwig Mies van der Rohe fought to keep the school open and moved the Bauhaus school to Berlin, but the attempt at survival was short-lived.  Suspected of producing anti-Nazi propaganda, the converted neglected factory was raided by the Gestapo.  Despite fighting to reopen the school for a time, Mies ultimately decided to close the school voluntarily, and he himself emigrated to the United States.  "
$insertRange.Font.Name = "Times New Roman"
$insertRange.Font.NameBi = "Times New Roman"
# paragraph 3
$insertRange.InsertParagraphAfter()
$insertRange = $d.Paragraphs.Item(15).Range
$insertRange.Text = "Mies was far from the first from the Bauhaus School to emigrate to the United States from Germany.  Many architects in particular “worked or tried to work for the National Socialist government in the years following Hitler’s ascent to power.  Only once it because clear that the Nazis were ruling modern architecture out of their agenda did these architects look for an exit” (Talesnik).  Those that preceded Mies included Walter Gropius, Marcel Breuer, Josef Albers, and László Moholy-Nagy, though this list omits plenty of strong emigrant contributors to the spread of Bauhaus.  Walter Gropius was announced the Chairman of the School of Architecture at Harvard in 1937, while others went on to be associated with other educational institutions including Yale and the Illinois Institute of Technology.  "
$insertRange.Font.Name = "Times New Roman"
$insertRange.Font.NameBi = "Times New Roman"
# paragraph 4
$insertRange.InsertParagraphAfter()
$insertRange = $d.Paragraphs.Item(16).Range
$insertRange.Text = "Irony lies in the execution of the National Socialist Party’s attempt to squash the Bauhaus school of thought and ideas.  In regarding the ‘foreignness’ of the school, they instead drove the free thinkers, both teachers and students of the school, abroad.  In doing so, they contributed to the spread of these ideas and spread the artistry and opportunities for manufacturing to other countries.  While many departed for Western Europe or the United States, those in later years would move “East to the Soviet Union and later to “other” countries like Chile, China, Hungary, Japan, Kenya, Mexico, North Korea, and Turkey”  (Talesnik).  The political pressure of the Socialist Party was a key driving factor in the global movement of those trained and practicing application of the teachings of the Bauhaus school, which in turn, helped the spread of Bauhaus itself.  This globalization meant that the Bauhaus principles were applied to different countries and markets, enabling worldwide adaptation that may otherwise have been centralized in Germany.  "
$insertRange.Font.Name = "Times New Roman"
$insertRange.Font.NameBi = "Times New Roman"
# paragraph 5
$insertRange.InsertParagraphAfter()
$insertRange = $d.Paragraphs.Item(17).Range
$insertRange.Text = "The influence of the Bauhaus movement has seen long and far-reaching impact on modern sensibilities, despite the school’s fairly modest 14 year history.  Yet the contributors, teachers, and students of the Bauhaus movement managed a legacy that still pervades our everyday lives.  While the school existed, the emphasis on practicality and affordable, simple design contributed to sustaining the school through profits of items sold.  The pressure and persecution of the Nazi Party against the “foreignness” of the Bauhaus drove the minds of the school abroad, ironically furthering the global adaptation of the Bauhaus movement.  It was a movement made for success, when practical, mass-produced artistry was able to reach many new markets, aided by the globalization forced by Nazi hands.  Embracing industrialization and commercialization meant that products were both affordable and beautiful; art made for everyday consumption.  These everyday objects have become entrenched in our everyday lives and families, to become our very own history.  Sale of such accessible items such as the Wassily Chair, or the influence of architecture on buildings that still stand erect today are simply evidence of the Bauhaus that has become entrenched in culture itself.  "
$insertRange.Font.Name = "Times New Roman"
$insertRange.Font.NameBi = "Times New Roman"
